$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-26 Monday" "2025-05-27 Tuesday"

Replace-Text "412×2=824" "679×3=2037"
Replace-Text "206×8=1648" "338×9=3042"
Replace-Text "324×9=2916" "309×2=618"
Replace-Text "171×4=684" "536×6=3216"
Replace-Text "915×6=5490" "199×6=1194"

Replace-Text "644×9=5796" "580×9=5220"
Replace-Text "643×7=4501" "725×9=6525"
Replace-Text "886×5=4430" "430×6=2580"
Replace-Text "752×7=5264" "460×7=3220"
Replace-Text "318×3=954" "767×9=6903"

Replace-Text "500×5=2500" "350×2=700"
Replace-Text "917×2=1834" "724×6=4344"
Replace-Text "681×6=4086" "894×3=2682"
Replace-Text "738×2=1476" "755×4=3020"
Replace-Text "548×4=2192" "659×7=4613"

Replace-Text "163×3=489" "226×8=1808"
Replace-Text "363×9=3267" "424×4=1696"
Replace-Text "880×6=5280" "294×6=1764"
Replace-Text "412×4=1648" "765×8=6120"
Replace-Text "961×2=1922" "758×7=5306"

Replace-Text "291×8=2328" "919×5=4595"
Replace-Text "473×7=3311" "167×4=668"
Replace-Text "873×2=1746" "435×3=1305"
Replace-Text "914×9=8226" "717×7=5019"
Replace-Text "217×9=1953" "232×2=464"
